$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.889.77"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "'1.732.54"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'245.18"
$ws.Range("E5").Value = "  +2.92%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "'0.4998"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("D8").Value = "'0.2711"
$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D9").Value = "'0.06150"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "'1.738.70"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").Value = "'0.07237"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("D12").Value = "'0.6520"
$ws.Range("E12").Value = "  +1.37%  "

$ws.Range("D13").Value = "'15.05"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "'4.756"
$ws.Range("E14").Value = "  +3.46%  "

$ws.Range("D15").Value = "'76.92"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "'25.901.18"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").Value = "'0.000006796"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'4.590"
$ws.Range("E21").Value = "  +7.32%  "

$ws.Range("D22").Value = "'1.961.47"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'8.748"
$ws.Range("E23").Value = "  +0.90%  "

$ws.Range("D24").Value = "'5.445"
$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("D25").Value = "'133.74"
$ws.Range("E25").Value = "  -3.62%  "

$ws.Range("D26").Value = "'15.22"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").Value = "'1.411"
$ws.Range("E27").Value = "  -7.70%  "

$ws.Range("D28").Value = "'1.775"
$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").Value = "'105.20"
$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("D30").Value = "'3.956"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("D31").Value = "'0.08093"
$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").Value = "'3.681"

$ws.Range("D33").Value = "'0.04732"
$ws.Range("E33").Value = "  +3.03%  "

$ws.Range("D34").Value = "'2.657"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "'0.9928"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").Value = "'0.6067"
$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("D37").Value = "'2.732"
$ws.Range("E37").Value = "  +1.62%  "

$ws.Range("D38").Value = "'0.01603"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.8552"
$ws.Range("E39").Value = "  +16.20%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.932"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").Value = "'99.93"
$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("D43").Value = "'0.3902"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").Value = "'5.009"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("D45").Value = "'0.1174"
$ws.Range("E45").Value = "  +4.15%  "

$ws.Range("D46").Value = "'6.311"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").Value = "'55.46"
$ws.Range("E47").Value = "  +0.80%  "

$ws.Range("D48").Value = "'0.05269"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").Value = "'30.63"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").Value = "'0.3470"
$ws.Range("E50").Value = "  +1.30%  "

$ws.Range("D51").Value = "'7.582"
$ws.Range("E51").Value = "  -0.71%  "
